$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column (C) for rows 2-9 from 45233 to 45243
foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = 45243
}
